$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.160.45"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.832.84"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6835"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3011"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07657"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").Value = "1.837.17"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.060"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6815"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.173"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.24%  "
$ws.Range("D17").Value = "29.163.49"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008177"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "2.081.71"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "226.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.418"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1455"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.743"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.510"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.256"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.140"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.201"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05154"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7672"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.837"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.131"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.675"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "1.307.99"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.724"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9389"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.798"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.982.99"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5200"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.533"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.770"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05919"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.96%  "
